$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("A4").Value = 7
$ws.Range("B6").Value = 2

# Move the mouse / selection from C8 to B8
$ws.Range("B8").Select()
